$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (Celaya vs Atlante). Row 3 (Leones Negros vs Tapatio) shifts up to become row 2.
$ws.Rows.Item(2).Delete()
